# Automatic tracker update: append new match rows (180-187) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row=180; EventId=14316289; Fecha="2025-08-13"; JugadorA="Frances Tiafoe";        JugadorB="Holger Rune";                 Pronostico="Gana Frances Tiafoe";              Cuota=2.1  },
    @{ Row=181; EventId=14316430; Fecha="2025-08-13"; JugadorA="Elena Rybakina";        JugadorB="Madison Keys";                Pronostico="Gana Madison Keys";                Cuota=2.2  },
    @{ Row=182; EventId=14369084; Fecha="2025-08-13"; JugadorA="Facundo Bagnis";        JugadorB="Bernard Tomic";               Pronostico="Gana Bernard Tomic";               Cuota=1.83 },
    @{ Row=183; EventId=14407242; Fecha="2025-08-13"; JugadorA="Andrea Pellegrino";     JugadorB="Coleman Wong";                Pronostico="Gana Andrea Pellegrino";           Cuota=4    },
    @{ Row=184; EventId=14407272; Fecha="2025-08-13"; JugadorA="Daniel Elahi Galan";    JugadorB="Ignacio Buse";                Pronostico="Gana Ignacio Buse";                Cuota=3.25 },
    @{ Row=185; EventId=14417899; Fecha="2025-08-13"; JugadorA="Liam Draxl";            JugadorB="Rio Noguchi";                 Pronostico="Gana Rio Noguchi";                 Cuota=3.75 },
    @{ Row=186; EventId=14407238; Fecha="2025-08-14"; JugadorA="Dusan Lajovic";         JugadorB="Roman Andres Burruchaga";     Pronostico="Gana Roman Andres Burruchaga";     Cuota=2.2  },
    @{ Row=187; EventId=14416065; Fecha="2025-08-13"; JugadorA="James McCabe";          JugadorB="Jacob Fearnley";              Pronostico="Gana James McCabe";                Cuota=3.5  }
)

foreach ($r in $newRows) {
    # A: event_id (number)
    $ws.Cells.Item($r.Row, 1).Value = $r.EventId

    # B: fecha -- keep it plain text ("2025-08-13"), not auto-converted to a date serial.
    $fechaCell = $ws.Cells.Item($r.Row, 2)
    $fechaCell.NumberFormat = "@"
    $fechaCell.Value = $r.Fecha
    $fechaCell.Style = "Normal"

    # C/D/E: player names + pronostico (plain text)
    $ws.Cells.Item($r.Row, 3).Value = $r.JugadorA
    $ws.Cells.Item($r.Row, 4).Value = $r.JugadorB
    $ws.Cells.Item($r.Row, 5).Value = $r.Pronostico

    # F: cuota (number)
    $ws.Cells.Item($r.Row, 6).Value = $r.Cuota

    # G/H: resultado / profit -- left blank, same as the other pending rows.
    $ws.Cells.Item($r.Row, 7).Style = "Normal"
    $ws.Cells.Item($r.Row, 8).Style = "Normal"
}
